$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Update the 2018 -> 2019 dollar-year references
$ws.Range("B26").Value = "2019 dollars per 2012 dollar"
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2019 dollar."'
$ws.Range("A21").Value = "million 2019 dollars"
$ws.Range("A18").Value = "billion 2019 dollars"

# Update the conversion-factor value
$ws.Range("A26").Value = 0.89805481563188172

$ws.Range("A19").Select()
